$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "英語" "English"
Replace-Text "葡萄牙語 / 法語 / 泰語 / 越南語 / 西班牙語" "Portuguese / French / Thai / Vietnamese / Spanish"
Replace-Text "簡介" "Brief"
Replace-Text "發送給在目標國家的合作夥伴的電子郵件，這些合作夥伴已回應參加，但在截止日期前未提交文件。 我們將取消他們的邀請。 將通過 customer.io 發送" "An email sent to partners in the target country who RSVPed yes but didn’t submit their documents by the deadline. We will be revoking their invites. It will be sent via customer.io"
Replace-Text "目標受眾" "Target audience"
Replace-Text "未按時提交文件的被邀請合作夥伴" "Invited partners who didn’t submit their documents on time"
Replace-Text "主題行" "Subject line"
Replace-Text ": 您的 " ": Your "
Replace-Text "[活動名稱]" "[EVENT NAME]"
Replace-Text " 註冊" " registration"
Replace-Text "沒有及時收到您的文件" "We didn’t receive your documents on time"
Replace-Text "[合作夥伴姓名]" "[PARTNER NAME]"
Replace-Text "截止日期（" "We didn’t receive your documents by the deadline ("
Replace-Text "[日月年]" "[DD Mmm YYYY]"
Replace-Text "）前沒有收到您的文件。 很遺憾，無法為您辦理 " "). Unfortunately, we’re unable to proceed with your registration for the "
Replace-Text "[活動名稱]" "[EVENT NAME]"
Replace-Text " 的註冊手續。" "."
Replace-Text "衷心祝愿您一切順利，並希望在下一次 " "We wish you the best and hope to see you at our next "
Replace-Text "會議/研討會/聯盟會員旅行" "conference/seminar/affiliate trip"
Replace-Text "中見到您。" "."

# Positional remap section (unique old texts).
# NOTE: order matters — the pre-existing "[NAME]" run must be renamed
# to "[WHATSAPP NO]" *before* a new "[NAME]" is created from
# "[電子郵件地址]", otherwise the later global Find would hit both.
Replace-Text "如有任何疑問，請通過 " "If you have any questions, please contact your country manager, "
Replace-Text "[WHATSAPP 號碼]" "[EMAIL ADDRESS]"
Replace-Text "[NAME]" "[WHATSAPP NO]"
Replace-Text "[電子郵件地址]" "[NAME]"
Replace-Text " 或 " ", at "
Replace-Text " (WhatsApp) 聯繫您的區域經理 " " or "
Replace-Text "。 " " (WhatsApp). "

# Comments
Replace-Text "選擇其中一個" "choose either one"
